$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Data": prepend 5 new weekly observations and append 2 new ones
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Data")

# Find the last populated row before any edits (should be row 113: 45252 / 7810.814)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# --- Insert 5 new rows right after the header row, shifting existing data down ---
$ws.Range("A2:A6").EntireRow.Insert()

# Copy the formatting (date number format for column A, plain number format for
# column B) from the row that used to be the first data row (now shifted down to
# row 7) onto the newly inserted rows.
$ws.Range("A7").Copy()
$ws.Range("A2:A6").PasteSpecial(-4122)
$ws.Range("B7").Copy()
$ws.Range("B2:B6").PasteSpecial(-4122)

$frontDates  = @(44440, 44447, 44454, 44461, 44468)
$frontValues = @(8349.173000000001, 8357.314, 8448.77, 8489.824000000001, 8447.981)
for ($i = 0; $i -lt $frontDates.Length; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $frontDates[$i]
    $ws.Cells.Item($r, 2).Value = $frontValues[$i]
}

# --- Append 2 new rows at the bottom ---
$newLastRow = $lastRow + 5
$ws.Range("A" + $newLastRow).Copy()
$ws.Range("A" + ($newLastRow + 1) + ":A" + ($newLastRow + 2)).PasteSpecial(-4122)
$ws.Range("B" + $newLastRow).Copy()
$ws.Range("B" + ($newLastRow + 1) + ":B" + ($newLastRow + 2)).PasteSpecial(-4122)

$endDates  = @(45259, 45266)
$endValues = @(7796.145, 7737.385)
for ($i = 0; $i -lt $endDates.Length; $i++) {
    $r = $newLastRow + 1 + $i
    $ws.Cells.Item($r, 1).Value = $endDates[$i]
    $ws.Cells.Item($r, 2).Value = $endValues[$i]
}

# ---------------------------------------------------------------------------
# Sheet "SeriesInfo": refresh the FRED metadata timestamps
# ---------------------------------------------------------------------------
$info = $wb.Worksheets.Item("SeriesInfo")

# These look like plain dates ("YYYY-MM-DD"), so force a text number format
# first to keep them stored as text instead of being auto-converted to date
# serial numbers.
$info.Range("B3").NumberFormat = "@"
$info.Range("B3").Value = "2023-12-08"
$info.Range("B4").NumberFormat = "@"
$info.Range("B4").Value = "2023-12-08"
$info.Range("B7").NumberFormat = "@"
$info.Range("B7").Value = "2023-12-06"

# This one includes a time + UTC offset, which Excel does not recognize as a
# date/time literal, so it is safely stored as text already.
$info.Range("B14").Value = "2023-12-07 15:34:03-06"

Write-Host "edit complete"
